$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume table with latest scraped values.
# Column D (Price) values are prefixed with a leading apostrophe so Excel
# keeps them as literal text (matching the source data) instead of
# auto-converting number-looking strings into numeric cells.
$ws.Range("D2").Value = "'46.159.05"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "'2.594.21"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'307.94"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "'99.03"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").Value = "'0.595"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.578"
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").Value = "'38.98"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'54.36"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "'0.0840"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "'8.14"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "'2.995.08"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "'2.598.32"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'0.916"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'14.84"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'46.281.67"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "'0.0000101"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "'12.83"
$ws.Range("E21").Value = "  -5.77%  "
$ws.Range("D22").Value = "'6.69"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'275.25"
$ws.Range("E23").Value = "  +8.65%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'71.42"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").Value = "'3.03"
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").Value = "'2.21"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").Value = "'29.61"
$ws.Range("E27").Value = "  +11.47%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'4.02"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("D30").Value = "'10.63"
$ws.Range("E30").Value = "  +2.38%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.20"
$ws.Range("E31").Value = "  -4.12%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'38.09"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").Value = "'6.30"
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("D34").Value = "'3.57"
$ws.Range("E34").Value = "  -8.25%  "
$ws.Range("D35").Value = "'2.22"
$ws.Range("E35").Value = "  -5.17%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'153.07"
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").Value = "'0.0834"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'2.78"
$ws.Range("E38").Value = "  -6.31%  "
$ws.Range("D39").Value = "'0.121"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").Value = "'23.13"
$ws.Range("E41").Value = "  +30.99%  "
$ws.Range("D42").Value = "'15.85"
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("D44").Value = "'3.57"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "'3.95"
$ws.Range("E45").Value = "  -5.57%  "
$ws.Range("D46").Value = "'2.109.13"
$ws.Range("E46").Value = "  +3.23%  "
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("D48").Value = "'95.44"
$ws.Range("E48").Value = "  +5.00%  "
$ws.Range("D49").Value = "'9.54"
$ws.Range("E49").Value = "  +6.48%  "
$ws.Range("D50").Value = "'108.43"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "'1.76"
$ws.Range("E51").Value = "  -2.21%  "
